$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Créditos-aula: 4 -> 2
# Leading apostrophe forces Excel to store this as text (matching the
# original shared-string cell type) instead of auto-converting to a number.
$ws.Range("B5").Value = "'2"
$ws.Range("C5").Value = "'2"

# Carga horária: 60 h -> 30 h
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Ativação: 01/01/2012 -> 01/01/2022
# Leading apostrophe forces Excel to store this as text (matching the
# original shared-string cell type) instead of auto-converting to a date serial.
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

# Objetivos: collapse embedded newlines into a single line
$objetivos = "Introduzir ao aluno a teoria de propriedades elétricas, térmicas, magnéticas e óticas de materiais sólidos, levando emconta o aspecto microscópico da estrutura do material. Dá-se ênfase à aplicação do material de acordo com aspropriedades que ele apresenta."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: collapse embedded newlines into a single line
$programaResumido = "PROPRIEDADES ELETRÔNICAS: Condutividade elétrica em metais, semicondutores e isolantes.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão, Condutividade e Tensões Térmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Lasers. Aplicações."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Programa: collapse embedded newlines into a single line
$programa = "PROPRIEDADES ELETRÔNICAS:Teoria do Elétron Livre em Metais. Níveis de Energia em Sólidos. Condutividade.Supercondutividade. Semicondutividade. Isolantes (Dielétricos). Aplicações.PROPRIEDADES MAGNÉTICAS: Conceitos Básicos. Curvas de Magnetização. Teoria de Domínio. MateriaisMagnéticos. Aplicações.PROPRIEDADES TÉRMICAS: Capacidade Calorífica. Expansão Térmica. Condutividade Térmica. TensõesTérmicas. AplicaçõesPROPRIEDADES ÓTICAS: Absorção. Transparência. Reflectividade. Fotocondutividade. Luminescência. Lasers.Fibra Ótica. Danos por Radiação. Aplicações."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa
